# Apply cryptocurrency price/volume updates (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.157.16'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '1.854.46'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("E4").Value = '  +0.10%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '235.98'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("E6").Value = '  +0.14%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4760'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -2.99%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2811'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -4.41%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06505'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.11%  '
$ws.Range("D10").Value = '1.862.01'
$ws.Range("E10").Value = '  -1.60%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07347'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.03%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '16.32'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -4.87%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.147'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.38%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '87.20'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.11%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.6445'
$cell.Style = "Normal"
$ws.Range("D16").Value = '30.121.84'
$ws.Range("E16").Value = '  -1.42%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '13.24'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("E18").Value = '  +0.09%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.000007611'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("D20").Value = '2.105.47'
$ws.Range("E20").Value = '  -1.73%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.01%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.252'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.87%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '217.73'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +14.21%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '6.104'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '9.292'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -1.96%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '165.83'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.29%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '18.54'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.07%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.907'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -1.53%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.424'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -3.60%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '4.248'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -3.05%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.09139'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -0.27%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.966'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -3.87%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.05027'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.80%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.7416'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.139'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +3.29%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.687'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -0.94%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.01820'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.609'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -3.00%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.9024'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -2.18%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '2.041'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.93%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '5.917'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '106.69'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.33%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.4245'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.89%  '
$ws.Range("E44").Value = '  +0.72%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '7.431'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.40%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.1310'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -5.38%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.561'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +9.45%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '64.08'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -7.51%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '8.814'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.57%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '34.25'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -2.05%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.05692'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -2.36%  '
